# Apply scheduled price/profit data refresh to each class sheet.
# Values below come from an external market-data pull; columns are:
# H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
# K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 616
$ws.Range("I41").Value = 138.5
$ws.Range("J41").Value = 1093.5
$ws.Range("K41").Value = 138.5
$ws.Range("L41").Value = 1093.5
$ws.Range("M41").Value = 301.5
$ws.Range("N41").Value = -1973.5
$ws.Range("H62").Value = 306690.5
$ws.Range("I62").Value = 457779.9
$ws.Range("K62").Value = 457779.9
$ws.Range("M62").Value = -457155.9
$ws.Range("H65").Value = 306690.5
$ws.Range("I65").Value = 457779.9
$ws.Range("K65").Value = 2288899.5
$ws.Range("M65").Value = -2285779.5
$ws.Range("H86").Value = 83337736
$ws.Range("I86").Value = 3829
$ws.Range("K86").Value = 3829
$ws.Range("M86").Value = -2706
$ws.Range("H89").Value = 83337736
$ws.Range("I89").Value = 3829
$ws.Range("K89").Value = 19145
$ws.Range("M89").Value = -13529
$ws.Range("H116").Value = 5100
$ws.Range("I116").Value = 5180
$ws.Range("K116").Value = 5180
$ws.Range("M116").Value = -1738
$ws.Range("H123").Value = 45894.445
$ws.Range("J123").Value = 45894.445
$ws.Range("L123").Value = 45894.445
$ws.Range("N123").Value = -55694.445

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1054.0869
$ws.Range("I45").Value = 976.26666
$ws.Range("K45").Value = 976.26666
$ws.Range("M45").Value = -599.26666
$ws.Range("H124").Value = 27283.3
$ws.Range("J124").Value = 27283.3
$ws.Range("L124").Value = 27283.3
$ws.Range("N124").Value = -37103.3
$ws.Range("H125").Value = 31632.637
$ws.Range("J125").Value = 31632.637
$ws.Range("L125").Value = 31632.637
$ws.Range("N125").Value = -41472.637

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H80").Value = 545.2353000000001
$ws.Range("J80").Value = 605.25
$ws.Range("L80").Value = 605.25
$ws.Range("N80").Value = -2601.25
$ws.Range("H83").Value = 545.2353000000001
$ws.Range("J83").Value = 605.25
$ws.Range("L83").Value = 3026.25
$ws.Range("N83").Value = -13010.25
$ws.Range("H94").Value = 822.26666
$ws.Range("I94").Value = 973.86957
$ws.Range("K94").Value = 973.86957
$ws.Range("M94").Value = -522.86957
$ws.Range("H113").Value = 23159.6
$ws.Range("I113").Value = 23159.6
$ws.Range("K113").Value = 23159.6
$ws.Range("M113").Value = -20989.6
$ws.Range("H134").Value = 1801.8889
$ws.Range("I134").Value = 1579.625
$ws.Range("J134").Value = 3580
$ws.Range("K134").Value = 4738.875
$ws.Range("L134").Value = 10740
$ws.Range("M134").Value = -2203.875
$ws.Range("N134").Value = -15810

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 900
$ws.Range("I10").Value = 900
$ws.Range("K10").Value = 900
$ws.Range("M10").Value = -761
$ws.Range("H31").Value = 39055.965
$ws.Range("I31").Value = 1863.7368
$ws.Range("J31").Value = 127387.5
$ws.Range("K31").Value = 1863.7368
$ws.Range("L31").Value = 127387.5
$ws.Range("M31").Value = -1568.7368
$ws.Range("N31").Value = -127977.5
$ws.Range("H34").Value = 39055.965
$ws.Range("I34").Value = 1863.7368
$ws.Range("J34").Value = 127387.5
$ws.Range("K34").Value = 1863.7368
$ws.Range("L34").Value = 127387.5
$ws.Range("M34").Value = -1661.7368
$ws.Range("N34").Value = -127791.5
$ws.Range("H62").Value = 5250
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376
$ws.Range("H65").Value = 5250
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 4566.6665
$ws.Range("I118").Value = 366.66666
$ws.Range("J118").Value = 6666.6665
$ws.Range("K118").Value = 1099.99998
$ws.Range("L118").Value = 19999.9995
$ws.Range("M118").Value = 143.0000199999999
$ws.Range("N118").Value = -22485.9995
$ws.Range("H127").Value = 1928.2858
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 1928.2858
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 5784.857400000001
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = -15704.8574
$ws.Range("H131").Value = 18594282
$ws.Range("J131").Value = 49267.81
$ws.Range("L131").Value = 147803.43
$ws.Range("N131").Value = -157883.43

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 800.6
$ws.Range("I3").Value = 750.75
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 750.75
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -634.75
$ws.Range("N3").Value = -1232
$ws.Range("H70").Value = 4088.625
$ws.Range("I70").Value = 3553.8462
$ws.Range("J70").Value = 4720.636
$ws.Range("K70").Value = 3553.8462
$ws.Range("L70").Value = 4720.636
$ws.Range("M70").Value = -3283.8462
$ws.Range("N70").Value = -5260.636
$ws.Range("H73").Value = 4088.625
$ws.Range("I73").Value = 3553.8462
$ws.Range("J73").Value = 4720.636
$ws.Range("K73").Value = 3553.8462
$ws.Range("L73").Value = 4720.636
$ws.Range("M73").Value = -2617.8462
$ws.Range("N73").Value = -6592.636
$ws.Range("H126").Value = 3040
$ws.Range("I126").Value = 3172.2222
$ws.Range("J126").Value = 1850
$ws.Range("K126").Value = 9516.6666
$ws.Range("L126").Value = 5550
$ws.Range("M126").Value = -7046.6666
$ws.Range("N126").Value = -10490

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4384.394
$ws.Range("I61").Value = 5010.926
$ws.Range("J61").Value = 1565
$ws.Range("K61").Value = 5010.926
$ws.Range("L61").Value = 1565
$ws.Range("M61").Value = -4808.926
$ws.Range("N61").Value = -1969
$ws.Range("H113").Value = 4384.394
$ws.Range("I113").Value = 5010.926
$ws.Range("J113").Value = 1565
$ws.Range("K113").Value = 5010.926
$ws.Range("L113").Value = 1565
$ws.Range("M113").Value = -2840.926
$ws.Range("N113").Value = -5905
$ws.Range("H122").Value = 13045
$ws.Range("I122").Value = 34666.668
$ws.Range("K122").Value = 104000.004
$ws.Range("M122").Value = -101550.004
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 43803.5
$ws.Range("I9").Value = 2600
$ws.Range("J9").Value = 85007
$ws.Range("K9").Value = 2600
$ws.Range("L9").Value = 85007
$ws.Range("M9").Value = -2460
$ws.Range("N9").Value = -85287
$ws.Range("H107").Value = 683.0454999999999
$ws.Range("I107").Value = 585.5263
$ws.Range("K107").Value = 1756.5789
$ws.Range("M107").Value = 163.4211
$ws.Range("H113").Value = 360.46155
$ws.Range("I113").Value = 381.2
$ws.Range("J113").Value = 347.5
$ws.Range("K113").Value = 1143.6
$ws.Range("L113").Value = 1042.5
$ws.Range("M113").Value = 1026.4
$ws.Range("N113").Value = -5382.5

Write-Output "Applied scheduled market-data refresh to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
